# Update the "想去人数" (want-to-go count) figures in the 展览 and 全部类型
# sheets to reflect the latest scrape, as described by the commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 260
    $ws.Range("F3").Value = 377
}
